$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.593.48'
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").Value = '1.671.27'
$ws.Range("E3").Value = '  -2.20%  '

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  +0.47%  '

$ws.Range("D5").Value = "'314.60"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").Value = "'0.3911"
$ws.Range("E7").Value = '  -3.16%  '

$ws.Range("D8").Value = "'0.3932"
$ws.Range("E8").Value = '  -3.55%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = "'51.96"
$ws.Range("E10").Value = '  -3.17%  '

$ws.Range("D11").Value = "'1.385"
$ws.Range("E11").Value = '  -6.48%  '

$ws.Range("D12").Value = "'0.08648"
$ws.Range("E12").Value = '  -2.15%  '

$ws.Range("D13").Value = "'25.06"
$ws.Range("E13").Value = '  -4.82%  '

$ws.Range("D14").Value = "'7.287"
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").Value = "'7.738"
$ws.Range("E15").Value = '  -4.97%  '

$ws.Range("D16").Value = "'0.00001309"
$ws.Range("E16").Value = '  -3.79%  '

$ws.Range("D17").Value = '1.677.76'
$ws.Range("E17").Value = '  -2.10%  '

$ws.Range("D18").Value = "'93.51"
$ws.Range("E18").Value = '  -3.64%  '

$ws.Range("D19").Value = "'0.07048"
$ws.Range("E19").Value = '  -1.60%  '

$ws.Range("D20").Value = "'20.56"
$ws.Range("E20").Value = '  -3.07%  '

$ws.Range("D21").Value = "'7.061"
$ws.Range("E21").Value = '  -3.35%  '

$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("D23").Value = "'13.94"
$ws.Range("E23").Value = '  -3.28%  '

$ws.Range("D24").Value = '24.631.62'
$ws.Range("E24").Value = '  -1.27%  '

$ws.Range("D25").Value = "'2.350"
$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("D26").Value = "'23.22"
$ws.Range("E26").Value = '  -0.35%  '

$ws.Range("D27").Value = "'2.714"
$ws.Range("E27").Value = '  -7.11%  '

$ws.Range("D28").Value = "'162.63"
$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("D29").Value = "'5.728"
$ws.Range("E29").Value = '  -8.13%  '

$ws.Range("D30").Value = "'147.12"
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").Value = "'7.840"
$ws.Range("E31").Value = '  -6.75%  '

$ws.Range("D32").Value = "'2.481"
$ws.Range("E32").Value = '  +11.01%  '

$ws.Range("D33").Value = '1.860.83'
$ws.Range("E33").Value = '  -2.65%  '

$ws.Range("D34").Value = "'0.08370"
$ws.Range("E34").Value = '  -6.06%  '

$ws.Range("D35").Value = "'0.03031"
$ws.Range("E35").Value = '  -5.61%  '

$ws.Range("D36").Value = "'6.896"
$ws.Range("E36").Value = '  -5.10%  '

$ws.Range("D37").Value = "'0.2799"
$ws.Range("E37").Value = '  -2.18%  '

$ws.Range("D38").Value = "'0.9765"
$ws.Range("E38").Value = '  -5.03%  '

$ws.Range("D39").Value = "'0.09453"
$ws.Range("E39").Value = '  +1.03%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'1.546"
$ws.Range("E40").Value = '  +4.99%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'10.51"
$ws.Range("E41").Value = '  -3.55%  '

$ws.Range("D42").Value = "'0.7906"
$ws.Range("E42").Value = '  -7.01%  '

$ws.Range("D43").Value = "'13.46"
$ws.Range("E43").Value = '  -5.29%  '

$ws.Range("D44").Value = "'16.36"
$ws.Range("E44").Value = '  -6.04%  '

$ws.Range("D45").Value = "'0.7100"
$ws.Range("E45").Value = '  -4.81%  '

$ws.Range("D46").Value = "'2.559"
$ws.Range("E46").Value = '  -5.96%  '

$ws.Range("D47").Value = "'4.208"
$ws.Range("E47").Value = '  -1.07%  '

$ws.Range("D48").Value = "'0.08626"
$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = '  +0.39%  '

$ws.Range("D50").Value = "'1.319"
$ws.Range("E50").Value = '  -5.85%  '

$ws.Range("D51").Value = "'136.90"
$ws.Range("E51").Value = '  -3.97%  '
